$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "Price" column (D) holds numeric-looking values stored as plain
# text (e.g. "95.389.52", "2.74") in the source XML. Temporarily force
# the column to Text format while writing so Excel does not silently
# convert these into numbers, then restore the original (default) style.
$ws.Range("D2:D50").NumberFormat = "@"

# Refresh Coin/Link/Price/Volume(1h) figures to the latest scrape.
# Rows 29 and 30 also swap content (WrappedeETH <-> Aptos) in this update.

$ws.Range("D2").Value = "95.389.52"
$ws.Range("E2").Value = "  -2.19%  "
$ws.Range("D3").Value = "3.614.86"
$ws.Range("E3").Value = "  -2.97%  "
$ws.Range("D4").Value = "2.74"
$ws.Range("E4").Value = "  +26.02%  "
$ws.Range("E5").Value = "  +0.01%  "
$ws.Range("D6").Value = "223.97"
$ws.Range("E6").Value = "  -5.98%  "
$ws.Range("D7").Value = "640.71"
$ws.Range("E7").Value = "  -2.51%  "
$ws.Range("D8").Value = "0.421"
$ws.Range("E8").Value = "  -5.12%  "
$ws.Range("E9").Value = "  +5.49%  "
$ws.Range("D10").Value = "0.999"
$ws.Range("E10").Value = "  -0.03%  "
$ws.Range("D11").Value = "3.606.76"
$ws.Range("E11").Value = "  -3.15%  "
$ws.Range("D12").Value = "50.47"
$ws.Range("E12").Value = "  +12.63%  "
$ws.Range("D13").Value = "0.217"
$ws.Range("E13").Value = "  +4.89%  "
$ws.Range("E14").Value = "  -6.78%  "
$ws.Range("D15").Value = "6.50"
$ws.Range("E15").Value = "  -5.12%  "
$ws.Range("D16").Value = "4.286.52"
$ws.Range("E16").Value = "  -3.06%  "
$ws.Range("D17").Value = "95.320.28"
$ws.Range("E17").Value = "  -2.02%  "
$ws.Range("D18").Value = "24.46"
$ws.Range("E18").Value = "  +29.52%  "
$ws.Range("D19").Value = "9.13"
$ws.Range("E19").Value = "  -1.41%  "
$ws.Range("D20").Value = "13.74"
$ws.Range("E20").Value = "  +4.88%  "
$ws.Range("D21").Value = "3.608.72"
$ws.Range("E21").Value = "  -3.09%  "
$ws.Range("D22").Value = "0.290"
$ws.Range("E22").Value = "  +36.47%  "
$ws.Range("D23").Value = "0.534"
$ws.Range("E23").Value = "  -1.39%  "
$ws.Range("D24").Value = "136.50"
$ws.Range("E24").Value = "  +16.26%  "
$ws.Range("D25").Value = "531.97"
$ws.Range("E25").Value = "  +0.56%  "
$ws.Range("E26").Value = "  -5.79%  "
$ws.Range("E27").Value = "  +1.48%  "
$ws.Range("E28").Value = "  -9.62%  "
$ws.Range("B29").Value = "Aptos"
$ws.Range("C29").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D29").Value = "13.20"
$ws.Range("E29").Value = "  -1.86%  "
$ws.Range("B30").Value = "WrappedeETH"
$ws.Range("C30").Value = "https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth"
$ws.Range("D30").Value = "3.777.64"
$ws.Range("E30").Value = "  -3.75%  "
$ws.Range("D31").Value = "13.35"
$ws.Range("E31").Value = "  +4.21%  "
$ws.Range("E32").Value = "  +3.31%  "
$ws.Range("E33").Value = "  +0.13%  "
$ws.Range("D34").Value = "1.88"
$ws.Range("E34").Value = "  +2.43%  "
$ws.Range("D35").Value = "0.640"
$ws.Range("E35").Value = "  +6.75%  "
$ws.Range("D36").Value = "33.72"
$ws.Range("E36").Value = "  +1.78%  "
$ws.Range("E38").Value = "  -0.06%  "
$ws.Range("D39").Value = "0.0554"
$ws.Range("E39").Value = "  +20.88%  "
$ws.Range("D41").Value = "8.57"
$ws.Range("E41").Value = "  -2.33%  "
$ws.Range("D42").Value = "7.31"
$ws.Range("E42").Value = "  +6.92%  "
$ws.Range("D43").Value = "590.93"
$ws.Range("E43").Value = "  -7.71%  "
$ws.Range("D44").Value = "0.504"
$ws.Range("E44").Value = "  +1.39%  "
$ws.Range("D45").Value = "1.02"
$ws.Range("E45").Value = "  +4.76%  "
$ws.Range("D46").Value = "40.95"
$ws.Range("E46").Value = "  +0.04%  "
$ws.Range("E47").Value = "  -0.58%  "
$ws.Range("E48").Value = "  -7.16%  "
$ws.Range("D49").Value = "9.32"
$ws.Range("E49").Value = "  +5.71%  "
$ws.Range("D50").Value = "234.58"
$ws.Range("E50").Value = "  +12.24%  "
$ws.Range("E51").Value = "  -2.45%  "

# Restore the default (General) style on the Price column now that the
# text values are safely written, matching the original formatting.
$ws.Range("D2:D50").Style = "Normal"
